$d = $word.ActiveDocument

# Locate the paragraph that contains the "LOM3003 ... (Requisito fraco)" text.
$found = $d.Content
$found.Find.Execute("LOM3003: Cinética de Transformação em Materiais (Requisito fraco)", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Determine the paragraph index of that match.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $found.Start -and $p.Range.End -ge $found.End) {
        $anchorIndex = $i
        break
    }
}

# Right after that paragraph there are four paragraphs to remove entirely:
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) an empty paragraph
#   4) an empty paragraph with a page break before it
$startPara = $d.Paragraphs.Item($anchorIndex + 1)
$endPara = $d.Paragraphs.Item($anchorIndex + 4)

$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$r.Delete()
